# Update target ratios (ion_ratio_avg / ion_ratio_cv) for the glucuronide
# compounds in the QA sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QA")

# Row 16: Morphine-3-B-D-glucuronide
$ws.Range("C16").Value = 12
$ws.Range("D16").Value = 0.5

# Row 17: Morphine-6-B-D-glucuronide
$ws.Range("C17").Value = 9.5500000000000007

# Row 18: Hydromorphone-3-B-D-glucuronide
$ws.Range("C18").Value = 2.67

# Row 19: Oxymorphone Glucuronide
$ws.Range("C19").Value = 4.32

# Row 20: Codeine-6-glucuronide
$ws.Range("C20").Value = 8.6199999999999992

# Row 21: Norbuprenorphine glucuronide
$ws.Range("C21").Value = 7.71

# Restore the active selection to where the author left off editing.
$ws.Range("C10").Select()
